# edit.ps1 - apply the changes described in the commit diff to the open
# presentation ($ppt.ActivePresentation):
#
#   1. Bump the cached "auto-update" date placeholder text on the slide
#      master and every slide layout from 12/9/2015 -> 12/10/2015 (the
#      field gets re-cached whenever the deck is re-saved on a later day).
#   2. Resize/reposition the "Gruppieren 30" group shape (the
#      SyncTool.FileSystem / SyncTool.FileSystem.Test box) on slide 1 so
#      its left edge moves left while its right edge stays put.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder text: "12/9/2015" -> "12/10/2015"
# ---------------------------------------------------------------------

function Update-DateText {
    param($shapes)

    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "12/9/2015") {
                $tr.Text = "12/10/2015"
            }
        }
    }
}

$master = $p.SlideMaster

# the slide master itself has one date placeholder
Update-DateText $master.Shapes

# ... and so does every custom (slide) layout under it
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-DateText $layout.Shapes
}

# ---------------------------------------------------------------------
# 2. Resize the "Gruppieren 30" group on slide 1
#    off  x: 4541661 -> 4466472   (left edge moves left by 75189 EMU)
#    ext cx: 6545176 -> 6620365   (width grows by 75189 EMU, right edge
#                                  stays fixed)
# ---------------------------------------------------------------------

$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "Gruppieren 30") {
        $shp.Left  = 351.6907196044922
        $shp.Width = 521.2886047363282
    }
}
